$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.176.28'
$ws.Range("E2").Value = '  -7.70%  '

$ws.Range("D3").Value = '3.648.03'
$ws.Range("E3").Value = '  -7.68%  '

$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.19%  '

$ws.Range("D5").Value = '574.39'
$ws.Range("E5").Value = '  -5.36%  '

$ws.Range("D6").Value = '168.65'
$ws.Range("E6").Value = '  -1.41%  '

$ws.Range("D7").Value = '3.636.50'
$ws.Range("E7").Value = '  -7.74%  '

$ws.Range("D8").Value = '0.615'
$ws.Range("E8").Value = '  -9.95%  '

$ws.Range("E9").Value = '  -0.01%  '

$ws.Range("D10").Value = '0.693'
$ws.Range("E10").Value = '  -11.98%  '

$ws.Range("D11").Value = '0.158'
$ws.Range("E11").Value = '  -12.42%  '

$ws.Range("D12").Value = '49.91'
$ws.Range("E12").Value = '  -11.02%  '

$ws.Range("D13").Value = '0.0000282'
$ws.Range("E13").Value = '  -13.64%  '

$ws.Range("D14").Value = '10.26'
$ws.Range("E14").Value = '  -10.94%  '

$ws.Range("D15").Value = '4.204.04'
$ws.Range("E15").Value = '  -8.30%  '

$ws.Range("D16").Value = '3.645.66'
$ws.Range("E16").Value = '  -7.93%  '

$ws.Range("E17").Value = '  -3.96%  '

$ws.Range("D18").Value = '19.08'
$ws.Range("E18").Value = '  -10.82%  '

$ws.Range("D19").Value = '12.63'
$ws.Range("E19").Value = '  -10.10%  '

$ws.Range("B20").Value = 'Polygon'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.10'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -10.78%  '

$ws.Range("B21").Value = 'WrappedBTC'
$ws.Range("C21").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D21").Value = '66.795.77'
$ws.Range("E21").Value = '  -8.14%  '

$ws.Range("D22").Value = '400.23'
$ws.Range("E22").Value = '  -10.11%  '

$ws.Range("D23").Value = '4.41'
$ws.Range("E23").Value = '  -9.24%  '

$ws.Range("D24").Value = '86.56'
$ws.Range("E24").Value = '  -9.60%  '

$ws.Range("D25").Value = '2.99'
$ws.Range("E25").Value = '  -10.92%  '

$ws.Range("D26").Value = '12.53'
$ws.Range("E26").Value = '  -11.97%  '

$ws.Range("D27").Value = '10.47'
$ws.Range("E27").Value = '  -7.49%  '

$ws.Range("E28").Value = '  +1.41%  '

$ws.Range("D29").Value = '3.73'
$ws.Range("E29").Value = '  -12.29%  '

$ws.Range("D30").Value = '9.29'
$ws.Range("E30").Value = '  -10.70%  '

$ws.Range("D31").Value = '32.02'
$ws.Range("E31").Value = '  -10.79%  '

$ws.Range("D32").Value = '7.33'
$ws.Range("E32").Value = '  -8.00%  '

$ws.Range("D33").Value = '12.21'
$ws.Range("E33").Value = '  -12.17%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '63.90'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -7.68%  '

$ws.Range("D35").Value = '0.114'
$ws.Range("E35").Value = '  -11.23%  '

$ws.Range("D36").Value = '42.43'
$ws.Range("E36").Value = '  -14.18%  '

$ws.Range("D37").Value = '585.24'
$ws.Range("E37").Value = '  -7.46%  '

$ws.Range("D38").Value = '0.0₃0876'
$ws.Range("E38").Value = '  -12.46%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +0.02%  '

$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  -0.15%  '

$ws.Range("D41").Value = '0.387'
$ws.Range("E41").Value = '  -10.09%  '

$ws.Range("D42").Value = '0.131'
$ws.Range("E42").Value = '  -10.01%  '

$ws.Range("D43").Value = '2.92'
$ws.Range("E43").Value = '  -15.76%  '

$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D44").Value = '2.59'
$ws.Range("E44").Value = '  -1.74%  '

$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = '0.0426'
$ws.Range("E45").Value = '  -11.00%  '

$ws.Range("D46").Value = '2.79'
$ws.Range("E46").Value = '  -12.01%  '

$ws.Range("D47").Value = '8.98'
$ws.Range("E47").Value = '  -15.14%  '

$ws.Range("D48").Value = '2.746.07'
$ws.Range("E48").Value = '  -2.92%  '

$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").Value = '0.132'
$ws.Range("E49").Value = '  -11.08%  '

$ws.Range("B50").Value = 'ApeXProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D50").Value = '3.11'
$ws.Range("E50").Value = '  -7.93%  '

$ws.Range("D51").Value = '2.63'
$ws.Range("E51").Value = '  -7.33%  '
